# Fruta / hortaliza, semanal
# The data rows 2..25 on the sheet get reshuffled: each destination row ends
# up holding the full contents (dates, volumes, prices, etc.) that used to
# live in a different source row. Columns A,B,C,E,F,G,I,N,O,Q,R never change
# value across this particular shuffle, but we copy the whole row (A..R) so
# the operation is robust regardless.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (rows identified by their original position)
$rowMap = @{
    2  = 9
    3  = 11
    4  = 21
    5  = 10
    6  = 20
    7  = 17
    8  = 5
    9  = 19
    10 = 2
    11 = 23
    12 = 15
    13 = 22
    14 = 3
    15 = 4
    16 = 25
    17 = 12
    18 = 24
    19 = 7
    20 = 8
    21 = 18
    22 = 16
    23 = 14
    24 = 13
    25 = 6
}

$firstCol = 1   # A
$lastCol  = 18  # R

# 1) Snapshot every source row's values before any writes happen, since the
#    mapping is a set of permutation cycles (writes to one row must not
#    clobber data another destination still needs to read).
$snapshot = @{}
for ($r = 2; $r -le 25; $r++) {
    $rowVals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# 2) Write each destination row from the snapshot of its mapped source row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
